$wb = $excel.ActiveWorkbook

$wsMeta = $wb.Worksheets.Item("Metadata")
$wsElem = $wb.Worksheets.Item("Elements")

# Version: 2.1.0 -> 2.2.0-ballot
$wsMeta.Range("B3").Value = "2.2.0-ballot"

# Date: 2025-12-18T17:25:31+00:00 -> 2025-12-19T08:32:44+00:00
$wsMeta.Range("B8").Value = "2025-12-19T08:32:44+00:00"

# Base Definition: add version suffix
$wsMeta.Range("B18").Value = "http://hl7.org/fhir/StructureDefinition/Extension|4.0.1"

# Binding Value Set: add version suffix
$wsElem.Range("Z6").Value = "https://interop.esante.gouv.fr/ig/fhir/tddui/ValueSet/tddui-discriminator-vs|2.2.0-ballot"

# Column Z (Binding Value Set) widened to fit the longer value (bestFit column)
$wsElem.Columns.Item(26).ColumnWidth = 67.78
